$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H86").Value = 192858480
$ws.Range("I86").Value = 116667300
$ws.Range("J86").Value = 250001870
$ws.Range("K86").Value = 116667300
$ws.Range("L86").Value = 250001870
$ws.Range("M86").Value = -116666177
$ws.Range("N86").Value = -250004116

$ws.Range("H87").Value = 100999.4
$ws.Range("I87").Value = 45000
$ws.Range("J87").Value = 114999.25
$ws.Range("K87").Value = 45000
$ws.Range("L87").Value = 114999.25
$ws.Range("M87").Value = -43752
$ws.Range("N87").Value = -117495.25

$ws.Range("H89").Value = 192858480
$ws.Range("I89").Value = 116667300
$ws.Range("J89").Value = 250001870
$ws.Range("K89").Value = 583336500
$ws.Range("L89").Value = 1250009350
$ws.Range("M89").Value = -583330884
$ws.Range("N89").Value = -1250020582

$ws.Range("H90").Value = 100999.4
$ws.Range("I90").Value = 45000
$ws.Range("J90").Value = 114999.25
$ws.Range("K90").Value = 135000
$ws.Range("L90").Value = 344997.75
$ws.Range("M90").Value = -128760
$ws.Range("N90").Value = -357477.75

$ws.Range("H92").Value = 875.8333
$ws.Range("I92").Value = 385.25
$ws.Range("J92").Value = 1857
$ws.Range("K92").Value = 385.25
$ws.Range("L92").Value = 1857
$ws.Range("M92").Value = 862.75
$ws.Range("N92").Value = -4353

$ws.Range("H113").Value = 3121.75
$ws.Range("I113").Value = 3121.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3121.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 132.25

$ws.Range("H125").Value = 3000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 27000
$ws.Range("N125").Value = -31920

$ws.Range("H137").Value = 485735.44
$ws.Range("I137").Value = 1598.2142
$ws.Range("J137").Value = 909355.5
$ws.Range("K137").Value = 4794.642599999999
$ws.Range("L137").Value = 2728066.5
$ws.Range("M137").Value = -2244.642599999999
$ws.Range("N137").Value = -2733166.5

$ws.Range("H138").Value = 83335090
$ws.Range("I138").Value = 1373
$ws.Range("J138").Value = 250002530
$ws.Range("K138").Value = 4119
$ws.Range("L138").Value = 750007590
$ws.Range("M138").Value = 1021
$ws.Range("N138").Value = -750017870

$ws.Range("H141").Value = 4985.091
$ws.Range("I141").Value = 3876.4443
$ws.Range("J141").Value = 9974
$ws.Range("K141").Value = 11629.3329
$ws.Range("L141").Value = 29922
$ws.Range("M141").Value = -6449.332900000001
$ws.Range("N141").Value = -40282

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1284.7727
$ws.Range("I2").Value = 913.25
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 913.25
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -800.25
$ws.Range("N2").Value = -5226

$ws.Range("H15").Value = 2687.5
$ws.Range("I15").Value = 850
$ws.Range("J15").Value = 3300
$ws.Range("K15").Value = 850
$ws.Range("L15").Value = 3300
$ws.Range("M15").Value = -500
$ws.Range("N15").Value = -4000

$ws.Range("H61").Value = 45468.87
$ws.Range("I61").Value = 2052.1177
$ws.Range("J61").Value = 168483
$ws.Range("K61").Value = 2052.1177
$ws.Range("L61").Value = 168483
$ws.Range("M61").Value = -1840.1177
$ws.Range("N61").Value = -168907

$ws.Range("H74").Value = 7019.1333
$ws.Range("I74").Value = 4293.3335
$ws.Range("J74").Value = 8836.333000000001
$ws.Range("K74").Value = 4293.3335
$ws.Range("L74").Value = 8836.333000000001
$ws.Range("M74").Value = -3419.3335
$ws.Range("N74").Value = -10584.333

$ws.Range("H77").Value = 7019.1333
$ws.Range("I77").Value = 4293.3335
$ws.Range("J77").Value = 8836.333000000001
$ws.Range("K77").Value = 21466.6675
$ws.Range("L77").Value = 44181.665
$ws.Range("M77").Value = -17098.6675
$ws.Range("N77").Value = -52917.665

$ws.Range("H110").Value = 1387.2727
$ws.Range("I110").Value = 973.3333
$ws.Range("J110").Value = 3250
$ws.Range("K110").Value = 973.3333
$ws.Range("L110").Value = 3250
$ws.Range("M110").Value = 1071.6667
$ws.Range("N110").Value = -7340

$ws.Range("H116").Value = 1284.7727
$ws.Range("I116").Value = 913.25
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 913.25
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 1380.75
$ws.Range("N116").Value = -9588

$ws.Range("H122").Value = 2796.4
$ws.Range("I122").Value = 2683
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 8049
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -5599
$ws.Range("N122").Value = -14650

$ws.Range("H136").Value = 45468.87
$ws.Range("I136").Value = 2052.1177
$ws.Range("J136").Value = 168483
$ws.Range("K136").Value = 6156.353099999999
$ws.Range("L136").Value = 505449
$ws.Range("M136").Value = -3606.353099999999
$ws.Range("N136").Value = -510549

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1284.7727
$ws.Range("I3").Value = 913.25
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 913.25
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -799.25
$ws.Range("N3").Value = -5228

$ws.Range("H64").Value = 2493.3333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2493.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 2493.3333
$ws.Range("N64").Value = -2943.3333
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 2493.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2493.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 2493.3333
$ws.Range("N67").Value = -4053.3333
$ws.Range("M67").ClearContents()

$ws.Range("H140").Value = 79995
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 79995
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 79995
$ws.Range("N140").Value = -90355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 3452.1667
$ws.Range("I7").Value = 67.59999999999999
$ws.Range("J7").Value = 20375
$ws.Range("K7").Value = 202.8
$ws.Range("L7").Value = 61125
$ws.Range("M7").Value = -90.79999999999998
$ws.Range("N7").Value = -61349

$ws.Range("H131").Value = 1661.5333
$ws.Range("I131").Value = 1104.1666
$ws.Range("J131").Value = 2033.1111
$ws.Range("K131").Value = 3312.4998
$ws.Range("L131").Value = 6099.3333
$ws.Range("M131").Value = 1727.5002
$ws.Range("N131").Value = -16179.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2184.25
$ws.Range("I80").Value = 1324.8
$ws.Range("J80").Value = 3616.6667
$ws.Range("K80").Value = 1324.8
$ws.Range("L80").Value = 3616.6667
$ws.Range("M80").Value = -326.8
$ws.Range("N80").Value = -5612.6667

$ws.Range("H83").Value = 2184.25
$ws.Range("I83").Value = 1324.8
$ws.Range("J83").Value = 3616.6667
$ws.Range("K83").Value = 6624
$ws.Range("L83").Value = 18083.3335
$ws.Range("M83").Value = -1632
$ws.Range("N83").Value = -28067.3335

$ws.Range("H108").Value = 88977.89
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 88977.89
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 88977.89
$ws.Range("N108").Value = -96657.89

$ws.Range("H109").Value = 86996.2
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 86996.2
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 86996.2
$ws.Range("N109").Value = -89076.2

$ws.Range("H110").Value = 99895.14
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 99895.14
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 99895.14
$ws.Range("N110").Value = -108075.14

$ws.Range("H132").Value = 8210.223
$ws.Range("I132").Value = 3752.6155
$ws.Range("J132").Value = 19800
$ws.Range("K132").Value = 11257.8465
$ws.Range("L132").Value = 59400
$ws.Range("M132").Value = -8727.8465
$ws.Range("N132").Value = -64460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6025.75
$ws.Range("I61").Value = 6172.2856
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 6172.2856
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -5970.2856
$ws.Range("N61").Value = -5404

$ws.Range("H109").Value = 92000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 92000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 92000
$ws.Range("N109").Value = -94774

$ws.Range("H113").Value = 6025.75
$ws.Range("I113").Value = 6172.2856
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 6172.2856
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -4002.2856
$ws.Range("N113").Value = -9340

$ws.Range("H136").Value = 4367.227
$ws.Range("I136").Value = 5470.3335
$ws.Range("J136").Value = 3043.5
$ws.Range("K136").Value = 16411.0005
$ws.Range("L136").Value = 9130.5
$ws.Range("M136").Value = -13861.0005
$ws.Range("N136").Value = -14230.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 75388
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 75388
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 75388
$ws.Range("N58").Value = -76004

$ws.Range("H122").Value = 3091.9375
$ws.Range("I122").Value = 3069.3572
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 9208.071599999999
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -6758.071599999999
$ws.Range("N122").Value = -14650

$ws.Range("H136").Value = 3074
$ws.Range("I136").Value = 2222
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 6666
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -4116
$ws.Range("N136").Value = -15600
